$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("AW2").Value = 204.083981
$ws.Range("AW3").Value = 47.199977
$ws.Range("AW4").Value = 34.958021
$ws.Range("AW5").Value = 159.029641
$ws.Range("AW6").Value = 161.816192
$ws.Range("AQ7").Value = 26.071806
$ws.Range("AW8").Value = 196.853692
$ws.Range("AW9").Value = 145.894329
$ws.Range("AW10").Value = 129.977616
$ws.Range("AW11").Value = 45.874259
$ws.Range("AW12").Value = 130.879988
$ws.Range("AW13").Value = 208.937685
$ws.Range("AW14").Value = 194.752836
$ws.Range("AW15").Value = 163.881366
$ws.Range("AW16").Value = 166.178669
$ws.Range("AW17").Value = 45.828484
$ws.Range("AW18").Value = 28.064745
$ws.Range("AW19").Value = 133.765498
$ws.Range("AW20").Value = 125.854259
$ws.Range("AW21").Value = 47.200405
$ws.Range("AQ22").Value = 26.071597
$ws.Range("AQ23").Value = 25.935961
$ws.Range("AW24").Value = 159.029109
$ws.Range("AW25").Value = 45.973322
$ws.Range("AT26").Value = 13.064699
$ws.Range("AW27").Value = 194.758345
$ws.Range("AW28").Value = 115.841829
$ws.Range("AW29").Value = 124.935313
$ws.Range("AW30").Value = 55.075359
$ws.Range("AW31").Value = 196.8486
$ws.Range("AW32").Value = 53.856759
$ws.Range("AW33").Value = 133.077315
$ws.Range("AW34").Value = 41.975104
$ws.Range("AW35").Value = 161.123113
$ws.Range("AW36").Value = 168.908854
$ws.Range("AW37").Value = 56.140752
$ws.Range("AK38").Value = 61.921343
$ws.Range("AW39").Value = 34.955
$ws.Range("AK40").Value = 61.921238
$ws.Range("AQ41").Value = 27.343275
$ws.Range("AN42").Value = 33.897951
$ws.Range("AW43").Value = 118.160116
$ws.Range("AW44").Value = 133.076701
$ws.Range("AW45").Value = 84.193472
$ws.Range("AW46").Value = 84.192257
$ws.Range("AK47").Value = 83.970856
$ws.Range("AW48").Value = 161.124097
$ws.Range("AW49").Value = 161.817604
$ws.Range("AW50").Value = 68.96022000000001
$ws.Range("AW51").Value = 175.03309
$ws.Range("AW52").Value = 152.790625
$ws.Range("AW53").Value = 199.895278
$ws.Range("AW54").Value = 189.944097
$ws.Range("AW55").Value = 20.046644
$ws.Range("AW56").Value = 161.816771
$ws.Range("AW57").Value = 56.144861
$ws.Range("AW58").Value = 35.764502
$ws.Range("AW59").Value = 97.81834499999999
$ws.Range("AW60").Value = 165.072454
$ws.Range("AW61").Value = 125.8736
$ws.Range("AW62").Value = 45.873472
$ws.Range("AW63").Value = 201.188414
$ws.Range("AW64").Value = 204.081308
$ws.Range("AW65").Value = 161.832546
$ws.Range("AW66").Value = 118.934074
$ws.Range("AW67").Value = 54.970671
$ws.Range("AW68").Value = 61.966632
$ws.Range("AW69").Value = 18.791875
$ws.Range("AQ70").Value = 25.094896
$ws.Range("AW71").Value = 14.094329
$ws.Range("AW72").Value = 175.903981
$ws.Range("AW73").Value = 112.930197
$ws.Range("AW74").Value = 47.857188
$ws.Range("AK75").Value = 18.948623
$ws.Range("AW76").Value = 117.138843
$ws.Range("AW77").Value = 34.945347
$ws.Range("AW78").Value = 189.947488
$ws.Range("AW79").Value = 84.193252
$ws.Range("AW80").Value = 69.108611
$ws.Range("AW81").Value = 63.803368
$ws.Range("AW82").Value = 53.855417
$ws.Range("AW83").Value = 168.972813
$ws.Range("AW84").Value = 199.847558
$ws.Range("AW85").Value = 147.817222
$ws.Range("AW86").Value = 117.13713
$ws.Range("AW87").Value = 97.81119200000001
$ws.Range("AW88").Value = 20.056007
$ws.Range("AW89").Value = 153.819734
$ws.Range("AW90").Value = 194.752454
$ws.Range("AW91").Value = 104.192882
$ws.Range("AW92").Value = 161.815995
$ws.Range("AT93").Value = 11.004502
$ws.Range("AW94").Value = 201.188808
$ws.Range("AW95").Value = 112.803472
$ws.Range("AW96").Value = 34.958576
$ws.Range("AW97").Value = 130.87265
$ws.Range("AW98").Value = 122.948507
$ws.Range("AW99").Value = 196.938414
$ws.Range("AW100").Value = 152.935405
$ws.Range("AW101").Value = 92.19415499999999
$ws.Range("AW102").Value = 33.973472
$ws.Range("AW103").Value = 151.37228
$ws.Range("AW104").Value = 160.015799
$ws.Range("AW105").Value = 112.859421
$ws.Range("AW106").Value = 61.948113
$ws.Range("AW107").Value = 34.944028
$ws.Range("AQ108").Value = 32.982222
$ws.Range("AT109").Value = 11.004074
$ws.Range("AW110").Value = 204.084271
$ws.Range("AW111").Value = 189.89515
$ws.Range("AW112").Value = 199.845694
$ws.Range("AW113").Value = 92.19255800000001
$ws.Range("AW114").Value = 63.794005
$ws.Range("AW115").Value = 45.972836
$ws.Range("AW116").Value = 178.814363
$ws.Range("AW117").Value = 152.955417
$ws.Range("AW118").Value = 130.825856
$ws.Range("AW119").Value = 117.134977
$ws.Range("AK120").Value = 61.921354
$ws.Range("AW121").Value = 194.753067
$ws.Range("AW122").Value = 34.945451
$ws.Range("AW123").Value = 18.792581
$ws.Range("AK124").Value = 18.948611
$ws.Range("AW125").Value = 189.847662
$ws.Range("AW126").Value = 20.04809
$ws.Range("AK127").Value = 61.921238
$ws.Range("AW128").Value = 104.055231
$ws.Range("AW129").Value = 98.206644
$ws.Range("AW130").Value = 207.897685
$ws.Range("AW131").Value = 133.738044
$ws.Range("AW132").Value = 92.192014
$ws.Range("AQ133").Value = 32.87044
$ws.Range("AW134").Value = 28.064514
$ws.Range("AW135").Value = 189.898356
$ws.Range("AW136").Value = 194.754722
$ws.Range("AW137").Value = 130.935347
$ws.Range("AK138").Value = 61.921273
$ws.Range("AW139").Value = 165.072025
$ws.Range("AW140").Value = 161.816493
$ws.Range("AW141").Value = 124.936343
$ws.Range("AQ142").Value = 15.189572
$ws.Range("AW143").Value = 97.06851899999999
$ws.Range("AW144").Value = 96.868576
$ws.Range("AW145").Value = 20.055637
$ws.Range("AW146").Value = 118.766713
$ws.Range("AW147").Value = 108.831944
$ws.Range("AW148").Value = 138.902118
$ws.Range("AW149").Value = 92.19131899999999
$ws.Range("AW150").Value = 35.765162
$ws.Range("AQ151").Value = 26.027118
$ws.Range("AW152").Value = 161.1239
$ws.Range("AW153").Value = 117.14059
$ws.Range("AW154").Value = 40.788113
$ws.Range("AW155").Value = 151.146713
$ws.Range("AW156").Value = 66.89315999999999
$ws.Range("AW157").Value = 34.942373
$ws.Range("AW158").Value = 199.848495
$ws.Range("AW159").Value = 176.198796
$ws.Range("AW160").Value = 115.84162
$ws.Range("AT161").Value = 11.737951
$ws.Range("AW162").Value = 151.151655
$ws.Range("AW163").Value = 204.084664
$ws.Range("AW164").Value = 125.852731
$ws.Range("AW165").Value = 196.848044
$ws.Range("AW166").Value = 187.954884
$ws.Range("AW167").Value = 125.873484
$ws.Range("AW168").Value = 54.970509
$ws.Range("AK169").Value = 83.970845
$ws.Range("AW170").Value = 196.847685
$ws.Range("AW171").Value = 152.951227
$ws.Range("AW172").Value = 98.20592600000001
$ws.Range("AW173").Value = 130.765637
$ws.Range("AW174").Value = 41.918981
$ws.Range("AW175").Value = 204.092234
$ws.Range("AW176").Value = 138.914306
$ws.Range("AW177").Value = 117.138102
$ws.Range("AW178").Value = 98.209468
$ws.Range("AW179").Value = 152.955046
$ws.Range("AK180").Value = 61.921354
$ws.Range("AW181").Value = 20.047407
$ws.Range("AW182").Value = 196.938692
$ws.Range("AW183").Value = 211.159525
$ws.Range("AW184").Value = 117.143947
$ws.Range("AW185").Value = 90.931343
$ws.Range("AW186").Value = 184.175231
$ws.Range("AW187").Value = 194.753634
$ws.Range("AW188").Value = 122.857188
$ws.Range("AW189").Value = 178.811204
$ws.Range("AW190").Value = 126.780625
$ws.Range("AW191").Value = 199.821829
$ws.Range("AW192").Value = 196.938148
$ws.Range("AW193").Value = 195.86037
$ws.Range("AW194").Value = 207.842951
$ws.Range("AW195").Value = 209.983345
$ws.Range("AW196").Value = 199.896528
$ws.Range("AW197").Value = 161.818218
$ws.Range("AW198").Value = 67.161227
$ws.Range("AW199").Value = 161.836192
$ws.Range("AK200").Value = 83.970856
$ws.Range("AW201").Value = 208.938669
$ws.Range("AW202").Value = 199.837593
$ws.Range("AW203").Value = 208.934537
$ws.Range("AW204").Value = 109.888843
$ws.Range("AW205").Value = 41.943137
$ws.Range("AW206").Value = 34.961285
$ws.Range("AQ207").Value = 26.071713
$ws.Range("AQ208").Value = 26.07191
$ws.Range("AW209").Value = 117.143299
$ws.Range("AW210").Value = 168.908773
$ws.Range("AW211").Value = 186.150081
$ws.Range("AW212").Value = 125.919688
$ws.Range("AW213").Value = 118.159884
$ws.Range("AW214").Value = 98.207778
$ws.Range("AW215").Value = 49.136644
$ws.Range("AW216").Value = 161.124433
$ws.Range("AW217").Value = 117.002894
$ws.Range("AW218").Value = 20.048634
$ws.Range("AW219").Value = 35.76566
$ws.Range("AW220").Value = 151.150428
$ws.Range("AW221").Value = 196.822928
$ws.Range("AW222").Value = 161.817002
$ws.Range("AW223").Value = 125.851076
$ws.Range("AW224").Value = 69.10566
$ws.Range("AW225").Value = 97.910544
$ws.Range("AW226").Value = 92.19044
$ws.Range("AK227").Value = 61.92125
$ws.Range("AW228").Value = 161.817407
$ws.Range("AW229").Value = 133.07794
$ws.Range("AN230").Value = 18.969201
$ws.Range("AW231").Value = 181.970602
$ws.Range("AQ232").Value = 32.946516
$ws.Range("AT233").Value = 11.773426
$ws.Range("AW234").Value = 161.085012
$ws.Range("AW235").Value = 92.200185
$ws.Range("AW236").Value = 92.198032
$ws.Range("AW237").Value = 66.78258099999999
$ws.Range("AW238").Value = 159.029271
$ws.Range("AW239").Value = 83.18875
$ws.Range("AW240").Value = 68.960972
$ws.Range("AW241").Value = 17.957662
$ws.Range("AW242").Value = 189.947708
$ws.Range("AW243").Value = 199.895787
$ws.Range("AW244").Value = 178.811551
$ws.Range("AW245").Value = 54.970405
$ws.Range("AW246").Value = 68.961782
$ws.Range("AW247").Value = 34.945174
$ws.Range("AW248").Value = 178.812095
$ws.Range("AW249").Value = 98.20217599999999
$ws.Range("AW250").Value = 49.949016
$ws.Range("AW251").Value = 188.840903
$ws.Range("AW252").Value = 104.056424
$ws.Range("AW253").Value = 34.961944
$ws.Range("AW254").Value = 56.143032
$ws.Range("AW255").Value = 175.032894
$ws.Range("AW256").Value = 55.076088
$ws.Range("AW257").Value = 45.872095
$ws.Range("AW258").Value = 28.904803
$ws.Range("AW259").Value = 84.192176
$ws.Range("AW260").Value = 176.199167
$ws.Range("AW261").Value = 187.955833
$ws.Range("AK262").Value = 18.948611
$ws.Range("AW263").Value = 151.151493
$ws.Range("AW264").Value = 130.87956
$ws.Range("AW265").Value = 72.206979
$ws.Range("AW266").Value = 49.942674
$ws.Range("AW267").Value = 175.161563
$ws.Range("AK268").Value = 83.970868
$ws.Range("AW269").Value = 189.160764
$ws.Range("AW270").Value = 133.768657
$ws.Range("AW271").Value = 194.751667
$ws.Range("AW272").Value = 117.144514
$ws.Range("AW273").Value = 161.81934
$ws.Range("AW274").Value = 34.813403
$ws.Range("AW275").Value = 186.150544
$ws.Range("AW276").Value = 117.137384
$ws.Range("AW277").Value = 91.948854
$ws.Range("AW278").Value = 14.094572
$ws.Range("AW279").Value = 126.910324
$ws.Range("AW280").Value = 201.190266
$ws.Range("AW281").Value = 90.862477
$ws.Range("AW282").Value = 45.972685
$ws.Range("AW283").Value = 20.06015
$ws.Range("AT284").Value = 11.00419
$ws.Range("AW285").Value = 90.861042
$ws.Range("AW286").Value = 90.862813
$ws.Range("AW287").Value = 194.752049
$ws.Range("AW288").Value = 147.813738
$ws.Range("AW289").Value = 66.894109
$ws.Range("AW290").Value = 147.816528
$ws.Range("AW291").Value = 161.121771
$ws.Range("AW292").Value = 204.085116
$ws.Range("AW293").Value = 98.201319
$ws.Range("AW294").Value = 61.9114
$ws.Range("AW295").Value = 49.946655
$ws.Range("AW296").Value = 182.859167
$ws.Range("AW297").Value = 137.918183
$ws.Range("AW298").Value = 20.061678
$ws.Range("AQ299").Value = 27.343148
$ws.Range("AW300").Value = 161.818669
$ws.Range("AW301").Value = 130.825556
$ws.Range("AW302").Value = 76.83627300000001
$ws.Range("AW303").Value = 56.139248
$ws.Range("AW304").Value = 45.872523
$ws.Range("AW305").Value = 140.843148
$ws.Range("AN306").Value = 20.779896
$ws.Range("AW307").Value = 133.074583
$ws.Range("AW308").Value = 66.89362300000001
$ws.Range("AW309").Value = 208.943495
$ws.Range("AW310").Value = 119.142465
$ws.Range("AW311").Value = 98.208681
$ws.Range("AW312").Value = 47.856701
$ws.Range("AW313").Value = 166.18184
$ws.Range("AQ314").Value = 32.831968
$ws.Range("AK315").Value = 61.921435
$ws.Range("AW316").Value = 28.905417
$ws.Range("AW317").Value = 122.902222
$ws.Range("AW318").Value = 101.840741
$ws.Range("AW319").Value = 66.894734
$ws.Range("AW320").Value = 47.857928
$ws.Range("AW321").Value = 175.161782
$ws.Range("AW322").Value = 104.055972
$ws.Range("AW323").Value = 92.18967600000001
$ws.Range("AW324").Value = 98.26218799999999
$ws.Range("AW325").Value = 33.976933
$ws.Range("AQ326").Value = 15.189722
$ws.Range("AW327").Value = 199.8464
$ws.Range("AW328").Value = 201.188576
$ws.Range("AK329").Value = 61.9214
$ws.Range("AW330").Value = 14.975868
$ws.Range("AW331").Value = 207.837245
$ws.Range("AW332").Value = 98.209896
$ws.Range("AW333").Value = 17.957581
$ws.Range("AW334").Value = 178.812616
$ws.Range("AW335").Value = 69.107164
$ws.Range("AW336").Value = 67.14219900000001
$ws.Range("AW337").Value = 194.7514
$ws.Range("AW338").Value = 166.179155
$ws.Range("AW339").Value = 84.779352
$ws.Range("AW340").Value = 168.971713
$ws.Range("AW341").Value = 117.138738
$ws.Range("AW342").Value = 203.977153
$ws.Range("AW343").Value = 199.883009
$ws.Range("AW344").Value = 47.863322
$ws.Range("AW345").Value = 201.189375
$ws.Range("AW346").Value = 178.813137
$ws.Range("AW347").Value = 208.934896
$ws.Range("AW348").Value = 56.871273
$ws.Range("AW349").Value = 194.749931
$ws.Range("AW350").Value = 161.81581
$ws.Range("AW351").Value = 152.948438
$ws.Range("AW352").Value = 69.106481
$ws.Range("AW353").Value = 20.060602
$ws.Range("AW354").Value = 14.093414
$ws.Range("AW355").Value = 115.841968
$ws.Range("AW356").Value = 115.841505
$ws.Range("AW357").Value = 20.059028
$ws.Range("AW358").Value = 34.959549
$ws.Range("AW359").Value = 47.85691
